# Generate Report for Handoff
# Adds two new localization entries (0098402b-... and 9a070ea2-...) as
# "Ready for handoff" rows to the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$ghBase   = "https://github.com/OpenLocalizationTest/oltest/blob/350fad43d11bee66309d339b3c9befbaef25f731/e2e"
$zhcnBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6359ef94d84201175cf78611e2eeda0b497eda96/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht"
$dedeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/233ab519a54be428b9cd7d48dca52ba176a36637/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht"

$id1 = "0098402b-a623-465f-9c4e-0baa758574cd"
$id2 = "9a070ea2-1528-4767-a3ae-27f6bb1bb311"

$zhfile1 = "$id1.b804eab168983cdb81cba4ec5a24a1217fd47928.zh-cn.xlf"
$zhfile2 = "$id2.1b637f98a0d62192cbf361667a19d981b2715cc9.zh-cn.xlf"
$defile1 = "$id1.b804eab168983cdb81cba4ec5a24a1217fd47928.de-de.xlf"
$defile2 = "$id2.1b637f98a0d62192cbf361667a19d981b2715cc9.de-de.xlf"

$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) -> matches the workbook's custom "HyperLink" style
$dateFormat     = "yyyy-mm-dd HH:mm:ss"

function Style-AsLink($range, $address, $display) {
    $range.Parent.Hyperlinks.Add($range, $address, "", "", $display) | Out-Null
    $range.Font.Color = $hyperlinkColor
    $range.Font.Underline = 2
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$id1.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-21 03:02:36"

$wsOverview.Range("A5").Value = "$id2.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-21 03:02:36"

Style-AsLink $wsOverview.Range("A4") "$ghBase/$id1.md" "$id1.md"
Style-AsLink $wsOverview.Range("A5") "$ghBase/$id2.md" "$id2.md"
Style-AsDate $wsOverview.Range("D4")
Style-AsDate $wsOverview.Range("D5")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = "$id1.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = $zhfile1
$wsZh.Range("E4").Value = "2016-03-21 03:02:27"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("J4").Value = "Include"

$wsZh.Range("A5").Value = "$id2.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = $zhfile2
$wsZh.Range("E5").Value = "2016-03-21 03:02:27"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("J5").Value = "Include"

Style-AsLink $wsZh.Range("A4") "$ghBase/$id1.md" "$id1.md"
Style-AsLink $wsZh.Range("D4") "$zhcnBase/$zhfile1" $zhfile1
Style-AsLink $wsZh.Range("A5") "$ghBase/$id2.md" "$id2.md"
Style-AsLink $wsZh.Range("D5") "$zhcnBase/$zhfile2" $zhfile2
Style-AsDate $wsZh.Range("E4")
Style-AsDate $wsZh.Range("H4")
Style-AsDate $wsZh.Range("E5")
Style-AsDate $wsZh.Range("H5")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = "$id1.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = $defile1
$wsDe.Range("E4").Value = "2016-03-21 03:02:36"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("J4").Value = "Include"

$wsDe.Range("A5").Value = "$id2.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = $defile2
$wsDe.Range("E5").Value = "2016-03-21 03:02:36"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("J5").Value = "Include"

Style-AsLink $wsDe.Range("A4") "$ghBase/$id1.md" "$id1.md"
Style-AsLink $wsDe.Range("D4") "$dedeBase/$defile1" $defile1
Style-AsLink $wsDe.Range("A5") "$ghBase/$id2.md" "$id2.md"
Style-AsLink $wsDe.Range("D5") "$dedeBase/$defile2" $defile2
Style-AsDate $wsDe.Range("E4")
Style-AsDate $wsDe.Range("H4")
Style-AsDate $wsDe.Range("E5")
Style-AsDate $wsDe.Range("H5")

Write-Host "Report rows for handoff added."
